# Edit: fix/complete the insurance (保險, sheet 5) and debt (債務, sheet 6)
# tables - add the missing metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index, ...) and
# correct the header rows that were incorrectly populated with data values
# instead of proper field-name labels.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-CellWithHeaderFormat {
    param($ws, $row, $col, $value)
    $ws.Cells.Item(1, 2).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-CellWithDataFormat {
    param($ws, $row, $col, $value)
    $ws.Cells.Item(2, 2).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($row, $col).Value = $value
}

# Sets a cell using the "data row" formatting, writing the value as plain
# text even when it looks like a date (e.g. "2013-12-26"), so Excel does
# not silently reinterpret it as a date serial number.
function Set-CellWithDataFormatAsText {
    param($ws, $row, $col, $value)
    $ws.Cells.Item(2, 2).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($row, $col).Value = "'" + $value
}

# ---------------------------------------------------------------------
# Sheet 5: 保險 (insurance)
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(5)

# Header row (row 1): columns B..K
Set-CellWithHeaderFormat $wsIns 1 2  "company"
Set-CellWithHeaderFormat $wsIns 1 3  "name"
Set-CellWithHeaderFormat $wsIns 1 4  "owner"
Set-CellWithHeaderFormat $wsIns 1 5  "property_category"
Set-CellWithHeaderFormat $wsIns 1 6  "category"
Set-CellWithHeaderFormat $wsIns 1 7  "date"
Set-CellWithHeaderFormat $wsIns 1 8  "legislator_name"
Set-CellWithHeaderFormat $wsIns 1 9  "legislator_id"
Set-CellWithHeaderFormat $wsIns 1 10 "source_file"
Set-CellWithHeaderFormat $wsIns 1 11 "index"

# Data rows 2..7
$insRows = @(
    @{ Row = 2; Index = 91; Company = "台灣人壽"; Name = "健康安心終身醫療B型"; Owner = "李桐豪" },
    @{ Row = 3; Index = 92; Company = "台灣人壽"; Name = "歲歲長泰還本終身";     Owner = "李桐豪" },
    @{ Row = 4; Index = 93; Company = "台灣人壽"; Name = "六六大順增額終身";     Owner = "黃素香" },
    @{ Row = 5; Index = 94; Company = "台灣人壽"; Name = "金寶貝兒童终身";       Owner = "黃素香" },
    @{ Row = 6; Index = 95; Company = "台灣人壽"; Name = "喜福還本定期";         Owner = "黃素香" },
    @{ Row = 7; Index = 96; Company = "台灣人壽"; Name = "金如意還本終身";       Owner = "黃素香" }
)

foreach ($r in $insRows) {
    $row = $r.Row
    Set-CellWithDataFormat       $wsIns $row 2  $r.Company
    Set-CellWithDataFormat       $wsIns $row 3  $r.Name
    Set-CellWithDataFormat       $wsIns $row 4  $r.Owner
    Set-CellWithDataFormat       $wsIns $row 5  "insurance"
    Set-CellWithDataFormat       $wsIns $row 6  "normal"
    Set-CellWithDataFormatAsText $wsIns $row 7  "2013-12-26"
    Set-CellWithDataFormat       $wsIns $row 8  "李桐豪"
    Set-CellWithDataFormat       $wsIns $row 9  896
    Set-CellWithDataFormat       $wsIns $row 10 "tmp2e9d1"
    Set-CellWithDataFormat       $wsIns $row 11 $r.Index
}

# ---------------------------------------------------------------------
# Sheet 6: 債務 (debt)
# ---------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(6)

# Header row (row 1): columns B..N
Set-CellWithHeaderFormat $wsDebt 1 2  "species"
Set-CellWithHeaderFormat $wsDebt 1 3  "debtor"
Set-CellWithHeaderFormat $wsDebt 1 4  "owner"
Set-CellWithHeaderFormat $wsDebt 1 5  "total"
Set-CellWithHeaderFormat $wsDebt 1 6  "register_date"
Set-CellWithHeaderFormat $wsDebt 1 7  "register_reason"
Set-CellWithHeaderFormat $wsDebt 1 8  "property_category"
Set-CellWithHeaderFormat $wsDebt 1 9  "category"
Set-CellWithHeaderFormat $wsDebt 1 10 "date"
Set-CellWithHeaderFormat $wsDebt 1 11 "legislator_name"
Set-CellWithHeaderFormat $wsDebt 1 12 "legislator_id"
Set-CellWithHeaderFormat $wsDebt 1 13 "source_file"
Set-CellWithHeaderFormat $wsDebt 1 14 "index"

# Data row 2
Set-CellWithDataFormat       $wsDebt 2 2  "抵押"
Set-CellWithDataFormat       $wsDebt 2 3  "李桐豪"
Set-CellWithDataFormat       $wsDebt 2 4  "台北富邦銀行臺北市中山區中山北路二段50號"
Set-CellWithDataFormat       $wsDebt 2 5  3322620
Set-CellWithDataFormatAsText $wsDebt 2 6  "102年03月11曰"
Set-CellWithDataFormat       $wsDebt 2 7  "購買房舍貸款"
Set-CellWithDataFormat       $wsDebt 2 8  "debt"
Set-CellWithDataFormat       $wsDebt 2 9  "normal"
Set-CellWithDataFormatAsText $wsDebt 2 10 "2013-12-26"
Set-CellWithDataFormat       $wsDebt 2 11 "李桐豪"
Set-CellWithDataFormat       $wsDebt 2 12 896
Set-CellWithDataFormat       $wsDebt 2 13 "tmp2e9d1"
Set-CellWithDataFormat       $wsDebt 2 14 106

Write-Host "Insurance and debt sheets updated."
